# Update cryptos list values/percentages per upstream refresh (GitHub Actions snapshot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.529.46"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.806.93"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'223.74"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'32.74"
$ws.Range("E8").Value = "  +3.41%  "
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "'0.0704"
$ws.Range("E10").Value = "  +6.34%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "2.067.37"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "'11.12"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "1.807.21"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "34.533.74"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "'4.30"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "'69.25"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "'251.60"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").Value = "0.0₃0797"
$ws.Range("E20").Value = "  +7.30%  "
$ws.Range("D21").Value = "'11.22"
$ws.Range("E21").Value = "  +6.51%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").Value = "'161.72"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'3.61"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("B34").Value = "Swop.fi"
$ws.Range("C34").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D34").Value = "'485.51"
$ws.Range("E34").Value = "  +831.37%  "
$ws.Range("D35").Value = "'1.90"
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").Value = "1.427.84"
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.650"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'0.970"
$ws.Range("E40").Value = "  +7.62%  "
$ws.Range("D41").Value = "'82.27"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("D45").Value = "'6.06"
$ws.Range("E45").Value = "  +3.25%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'12.51"
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("B48").Value = "Kaspa"
$ws.Range("C48").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D48").Value = "'0.0497"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "1.959.22"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "'105.67"
$ws.Range("E50").Value = "  +7.68%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.04%  "

Write-Host "Updated cryptos sheet with latest snapshot values"
